$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the trial-5 row values (row 6: y_nrSteps, alienID, praclen)
$ws.Range("E6").Value = 7
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 13

# Move the active selection to match the saved cursor position
$ws.Range("E6").Select()
